$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: DATE_TYPE_CODE "001" -> "002" (force text so it doesn't become numeric 2)
$ws.Range("J2").Value = "'002"
$ws.Range("J2").Style = "Normal"

# N2: REPORT_DATE text update
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# O2..Y2 numeric updates
$ws.Range("O2").Value = 543932526.79
$ws.Range("P2").Value = 1793011.89
$ws.Range("Q2").Value = 170125991.84
$ws.Range("R2").Value = 35.2981385961
$ws.Range("S2").Value = 157820711.15
$ws.Range("T2").Value = 38.9306190033
$ws.Range("U2").Value = 30560906.47
$ws.Range("V2").Value = -16.799397901
$ws.Range("W2").Value = 161844007.91
$ws.Range("X2").Value = 19036403.78
$ws.Range("Y2").Value = 66.1520397362

# Z2 / AA2 cleared to empty text cells (still present, but blank)
$ws.Range("Z2").Value = "'"
$ws.Range("AA2").Value = "'"
$ws.Range("Z2").Style = "Normal"
$ws.Range("AA2").Style = "Normal"

# AB2..AG2 numeric updates
$ws.Range("AB2").Value = 382088518.88
$ws.Range("AC2").Value = 28.5415024069
$ws.Range("AD2").Value = 45.1796863483
$ws.Range("AE2").Value = 109.067163767
$ws.Range("AF2").Value = 418.1649147499
$ws.Range("AG2").Value = 29.7544272385
